$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, taken from the target diff.
$changes = @(
    @{cell='D2'; val='23.824.50'},
    @{cell='E2'; val='  -3.12%  '},
    @{cell='D3'; val='1.616.95'},
    @{cell='E3'; val='  -3.52%  '},
    @{cell='D4'; val='1.001'},
    @{cell='E4'; val='  -0.14%  '},
    @{cell='D5'; val='306.89'},
    @{cell='E5'; val='  -2.10%  '},
    @{cell='E6'; val='  -0.06%  '},
    @{cell='D7'; val='0.3915'},
    @{cell='E7'; val='  -0.32%  '},
    @{cell='D8'; val='0.3824'},
    @{cell='E8'; val='  -3.36%  '},
    @{cell='D9'; val='1.001'},
    @{cell='E9'; val='  -0.16%  '},
    @{cell='D10'; val='1.360'},
    @{cell='E10'; val='  -3.14%  '},
    @{cell='D11'; val='49.12'},
    @{cell='E11'; val='  -3.35%  '},
    @{cell='D12'; val='0.08397'},
    @{cell='E12'; val='  -3.10%  '},
    @{cell='D13'; val='23.93'},
    @{cell='E13'; val='  -5.51%  '},
    @{cell='D14'; val='7.024'},
    @{cell='E14'; val='  -4.33%  '},
    @{cell='D15'; val='7.536'},
    @{cell='E15'; val='  -2.46%  '},
    @{cell='D16'; val='0.00001274'},
    @{cell='E16'; val='  -3.55%  '},
    @{cell='D17'; val='1.601.93'},
    @{cell='E17'; val='  -4.03%  '},
    @{cell='D18'; val='93.04'},
    @{cell='E18'; val='  -1.04%  '},
    @{cell='D19'; val='0.06896'},
    @{cell='E19'; val='  -1.83%  '},
    @{cell='D20'; val='20.06'},
    @{cell='E20'; val='  -5.44%  '},
    @{cell='D21'; val='6.819'},
    @{cell='D22'; val='1.0000'},
    @{cell='E22'; val='  -0.02%  '},
    @{cell='D23'; val='13.40'},
    @{cell='E23'; val='  -4.14%  '},
    @{cell='D24'; val='23.820.75'},
    @{cell='D25'; val='2.445'},
    @{cell='E25'; val='  +3.32%  '},
    @{cell='D26'; val='2.847'},
    @{cell='E26'; val='  +2.24%  '},
    @{cell='D27'; val='22.13'},
    @{cell='E27'; val='  -4.41%  '},
    @{cell='D28'; val='157.08'},
    @{cell='E28'; val='  -2.03%  '},
    @{cell='D29'; val='139.19'},
    @{cell='E29'; val='  -5.10%  '},
    @{cell='D30'; val='5.244'},
    @{cell='E30'; val='  -10.70%  '},
    @{cell='D31'; val='7.837'},
    @{cell='E31'; val='  -5.74%  '},
    @{cell='D32'; val='2.484'},
    @{cell='E32'; val='  -1.04%  '},
    @{cell='D33'; val='1.796.44'},
    @{cell='E33'; val='  -2.39%  '},
    @{cell='D34'; val='0.08067'},
    @{cell='E34'; val='  -3.01%  '},
    @{cell='D35'; val='0.9758'},
    @{cell='E35'; val='  -1.31%  '},
    @{cell='E36'; val='  -7.22%  '},
    @{cell='D37'; val='6.581'},
    @{cell='E37'; val='  -5.59%  '},
    @{cell='D38'; val='0.2659'},
    @{cell='E38'; val='  -5.45%  '},
    @{cell='D39'; val='0.09186'},
    @{cell='E39'; val='  -3.68%  '},
    @{cell='D40'; val='10.28'},
    @{cell='E40'; val='  -0.74%  '},
    @{cell='D41'; val='13.41'},
    @{cell='E41'; val='  -1.26%  '},
    @{cell='D42'; val='1.430'},
    @{cell='E42'; val='  -5.44%  '},
    @{cell='D43'; val='0.7463'},
    @{cell='E43'; val='  -5.73%  '},
    @{cell='D44'; val='15.98'},
    @{cell='E44'; val='  -3.83%  '},
    @{cell='D45'; val='0.6851'},
    @{cell='E45'; val='  -3.89%  '},
    @{cell='D46'; val='2.455'},
    @{cell='E46'; val='  -4.27%  '},
    @{cell='E47'; val='  -2.67%  '},
    @{cell='E48'; val='  +0.00%  '},
    @{cell='D49'; val='0.08262'},
    @{cell='E49'; val='  -4.36%  '},
    @{cell='D50'; val='133.10'},
    @{cell='E50'; val='  -3.35%  '},
    @{cell='D51'; val='1.208'},
    @{cell='E51'; val='  -9.33%  '}
)

# All D/E values in this sheet are plain text (coinranking price/volume strings),
# even when they look like numbers (e.g. "1.001", "0.00001274"). Assigning such a
# string straight to Range.Value lets Excel auto-coerce it into a real number, which
# would not match the source data. Routing the text through a helper cell's formula
# ( ="literal text" ) and pasting *values only* keeps the destination a plain text
# cell/string with no numeric coercion and without touching any cell styles.
$scratch = $ws.Range("Z100")
foreach ($item in $changes) {
    $scratch.Formula = '="' + $item.val + '"'
    $scratch.Copy()
    $ws.Range($item.cell).PasteSpecial(-4163)
}
$scratch.Clear()
$excel.CutCopyMode = $false
Write-Output "Applied $($changes.Count) cell updates"
